$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Collapse every column on the sheet (columns A:K carry the explicit
#     custom widths, column L stands in for the trailing 12..16384 default
#     span) so the whole column range is grouped/collapsed, matching the
#     "collapsed" toggle applied uniformly across every <col> definition. ---
$ws.Columns("A:K").Group()
$ws.Columns("L").Group()
$ws.Outline.ShowLevels(0, 0)

# --- D10 changes from 21 to 100.0, mirroring C10's numeric value/style. ---
$ws.Range("D10").Value = 100.0
